$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row pairs whose B,C,D,E,F,G values need to be swapped (per commit diff)
$pairs = @(
    @(162, 163),
    @(183, 184),
    @(279, 280),
    @(317, 318),
    @(346, 347),
    @(351, 352),
    @(355, 356),
    @(372, 373),
    @(379, 380),
    @(400, 401),
    @(431, 432),
    @(457, 458),
    @(579, 580),
    @(581, 582),
    @(583, 584),
    @(586, 587),
    @(593, 594),
    @(601, 602),
    @(715, 716),
    @(720, 721)
)

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]
    for ($col = 2; $col -le 7; $col++) {
        $v1 = $ws.Cells.Item($r1, $col).Value2
        $v2 = $ws.Cells.Item($r2, $col).Value2
        $ws.Cells.Item($r1, $col).Value = $v2
        $ws.Cells.Item($r2, $col).Value = $v1
    }
}
